$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.03578230111377778
$ws.Range("R2").Value = 0.322040710024
$ws.Range("S2").Value = 0.003499619873322347
$ws.Range("T2").Value = 0.003499619873322347

# Row 3 updates
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("S3").Value = 0.8692174743460166
$ws.Range("T3").Value = 0.8692174743460165

# Row 4 updates
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("S4").Value = 0.1272829057806611
$ws.Range("T4").Value = 0.1272829057806611
